$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update timestamp header (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 13 de Julio de 2020 a las 19:48"

# --- Country name re-sorts (shared-string swaps realised as cell text updates) ---
$ws.Range("A128").Value = "Libia"
$ws.Range("A129").Value = "Yemen"
$ws.Range("A130").Value = "Suazilandia"
$ws.Range("A131").Value = "Benin"
$ws.Range("A135").Value = "Mozambique"
$ws.Range("A136").Value = "Jordania"
$ws.Range("A137").Value = "Letonia"
$ws.Range("A164").Value = "Birmania"
$ws.Range("A165").Value = "Isla de Man"
$ws.Range("A209").Value = "Groenlandia"
$ws.Range("A210").Value = "Islas Malvinas"

# --- Updated case numbers (daily data refresh) ---
$ws.Range("B4").Value = 3441503
$ws.Range("C4").Value = 27508
$ws.Range("D4").Value = 1534925
$ws.Range("E4").Value = 1768631
$ws.Range("G4").Value = 165
$ws.Range("H4").Value = 137947
$ws.Range("B6").Value = 906617
$ws.Range("C6").Value = 27151
$ws.Range("D6").Value = 571578
$ws.Range("E6").Value = 311312
$ws.Range("G6").Value = 540
$ws.Range("H6").Value = 23727
$ws.Range("B18").Value = 214001
$ws.Range("C18").Value = 1008
$ws.Range("D18").Value = 195671
$ws.Range("E18").Value = 12948
$ws.Range("G18").Value = 19
$ws.Range("H18").Value = 5382
$ws.Range("B24").Value = 104016
$ws.Range("C24").Value = 418
$ws.Range("D24").Value = 100627
$ws.Range("E24").Value = 3240
$ws.Range("G24").Value = 2
$ws.Range("H24").Value = 149
$ws.Range("B46").Value = 40248
$ws.Range("C46").Value = 1578
$ws.Range("D46").Value = 19323
$ws.Range("E46").Value = 20560
$ws.Range("G46").Value = 3
$ws.Range("H46").Value = 365
$ws.Range("B56").Value = 25638
$ws.Range("C56").Value = 10
$ws.Range("E56").Value = 528
$ws.Range("B65").Value = 15936
$ws.Range("C65").Value = 191
$ws.Range("D65").Value = 12934
$ws.Range("E65").Value = 2747
$ws.Range("G65").Value = 5
$ws.Range("H65").Value = 255
$ws.Range("B108").Value = 2762
$ws.Range("C108").Value = 31
$ws.Range("D108").Value = 2290
$ws.Range("E108").Value = 459
$ws.Range("B125").Value = 1642
$ws.Range("C125").Value = 7
$ws.Range("D125").Value = 1175
$ws.Range("E125").Value = 404
$ws.Range("E127").Value = 297
$ws.Range("G127").Value = 1
$ws.Range("H127").Value = 8
$ws.Range("B128").Value = 1512
$ws.Range("C128").Value = 79
$ws.Range("D128").Value = 367
$ws.Range("E128").Value = 1105
$ws.Range("G128").Value = 1
$ws.Range("H128").Value = 40
$ws.Range("B129").Value = 1498
$ws.Range("C129").Value = 33
$ws.Range("D129").Value = 675
$ws.Range("E129").Value = 399
$ws.Range("G129").Value = 7
$ws.Range("H129").Value = 424
$ws.Range("B130").Value = 1389
$ws.Range("C130").Value = 38
$ws.Range("D130").Value = 688
$ws.Range("E130").Value = 681
$ws.Range("H130").Value = 20
$ws.Range("B131").Value = 1378
$ws.Range("D131").Value = 557
$ws.Range("E131").Value = 795
$ws.Range("H131").Value = 26
$ws.Range("B135").Value = 1219
$ws.Range("C135").Value = 62
$ws.Range("D135").Value = 369
$ws.Range("E135").Value = 841
$ws.Range("H135").Value = 9
$ws.Range("B136").Value = 1183
$ws.Range("C136").Value = 4
$ws.Range("D136").Value = 1008
$ws.Range("E136").Value = 165
$ws.Range("G136").Value = 0
$ws.Range("H136").Value = 10
$ws.Range("B137").Value = 1174
$ws.Range("C137").Value = 1
$ws.Range("D137").Value = 1019
$ws.Range("E137").Value = 124
$ws.Range("G137").Value = 1
$ws.Range("H137").Value = 31
$ws.Range("C164").Value = 5
$ws.Range("D164").Value = 261
$ws.Range("E164").Value = 69
$ws.Range("H164").Value = 6
$ws.Range("B165").Value = 336
$ws.Range("D165").Value = 312
$ws.Range("E165").Value = 0
$ws.Range("H165").Value = 24
